$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.368.50'
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").Value = '3.677.30'
$ws.Range("E3").Value = '  -0.34%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '684.62'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.90'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -2.24%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  -1.32%  '

$ws.Range("E9").Value = '  -1.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.03'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -3.46%  '

$ws.Range("E11").Value = '  -3.31%  '

$ws.Range("E12").Value = '  -1.34%  '

$ws.Range("D13").Value = '4.299.53'
$ws.Range("E13").Value = '  -0.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.28'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -3.69%  '

$ws.Range("D15").Value = '3.677.31'
$ws.Range("E15").Value = '  -0.38%  '

$ws.Range("D16").Value = '69.339.23'

$ws.Range("E17").Value = '  +2.27%  '

$ws.Range("E18").Value = '  -3.14%  '

$ws.Range("E19").Value = '  -3.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '469.42'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -2.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.95'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +1.49%  '

$ws.Range("E22").Value = '  -2.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.87'
$ws.Range("D23").NumberFormat = "General"

$ws.Range("D24").Value = '3.823.68'
$ws.Range("E24").Value = '  -0.34%  '

$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("E26").Value = '  -4.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.93'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -5.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.18'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -4.16%  '

$ws.Range("E29").Value = '  -1.36%  '

$ws.Range("E30").Value = '  -5.54%  '

$ws.Range("E31").Value = '  +0.11%  '

$ws.Range("E32").Value = '  -4.18%  '

$ws.Range("E33").Value = '  -6.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.89'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -0.65%  '

$ws.Range("D35").Value = '3.652.94'
$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.159'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -4.65%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.17'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -4.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.15'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +1.98%  '

$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("E40").Value = '  +2.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0900'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -4.48%  '

$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.940'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -2.33%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '166.33'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +5.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.49'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -1.48%  '

$ws.Range("E46").Value = '  +1.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.70'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -4.75%  '

$ws.Range("E48").Value = '  +4.47%  '

$ws.Range("E49").Value = '  +0.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.78'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -3.78%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '27.12'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -1.75%  '

